$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $excel.Width
Write-Host $excel.Left
try { Write-Host $excel.ActiveWindow.Width } catch { Write-Host "no ActiveWindow" }
